$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "other_notes"
$ws.Range("I25").Value = "Data are projections for 2040-2059"
$ws.Range("I26").Value = "Data are projections for 2040-2059."
$ws.Range("I27").Value = "Data are projections for 2040-2059"

$note = "Data retrieved via API in March 2019. For detailed information on the observation level (e.g. National Estimation, UIS Estimation, or Category not applicable), please visit UIS.Stat (http://data.uis.unesco.org/)."
$ws.Range("I53").Value = $note
$ws.Range("I56").Value = $note
$ws.Range("I57").Value = $note
$ws.Range("I58").Value = $note
$ws.Range("I59").Value = $note
$ws.Range("I60").Value = $note
$ws.Range("I107").Value = $note
$ws.Range("I114").Value = $note
$ws.Range("I119").Value = $note
$ws.Range("I132").Value = $note
$ws.Range("I135").Value = $note
$ws.Range("I136").Value = $note

$ws.Range("I1").Select()
